$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'48.298.26"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.84%  "
$ws.Range("D3").Value = "'2.505.85"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.62%  "
$ws.Range("D5").Value = "'321.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("D6").Value = "'108.10"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.94%  "
$ws.Range("E7").Value = "  +1.40%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "'0.541"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("D10").Value = "'39.87"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.23%  "
$ws.Range("D11").Value = "'20.27"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +8.84%  "
$ws.Range("D12").Value = "'0.0820"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.16%  "
$ws.Range("E13").Value = "  -0.19%  "
$ws.Range("D14").Value = "'7.19"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.08%  "
$ws.Range("D15").Value = "'2.897.79"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.69%  "
$ws.Range("D16").Value = "'2.505.41"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.69%  "
$ws.Range("D17").Value = "'0.845"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.19%  "
$ws.Range("D18").Value = "'48.142.44"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.80%  "
$ws.Range("D19").Value = "'13.10"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.62%  "
$ws.Range("E20").Value = "  +1.68%  "
$ws.Range("D21").Value = "'0.0₃0945"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.30%  "
$ws.Range("E22").Value = "  +1.04%  "
$ws.Range("D23").Value = "'280.77"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +13.70%  "
$ws.Range("D24").Value = "'72.32"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.36%  "
$ws.Range("D25").Value = "'2.56"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.53%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").Value = "'25.77"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.06%  "
$ws.Range("D28").Value = "'2.30"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("E29").Value = "  -1.94%  "
$ws.Range("E30").Value = "  +0.95%  "
$ws.Range("D31").Value = "'35.27"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.58%  "
$ws.Range("E32").Value = "  -1.07%  "
$ws.Range("D33").Value = "'19.65"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.91%  "
$ws.Range("E34").Value = "  +0.90%  "
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("E36").Value = "  -0.51%  "
$ws.Range("E37").Value = "  -0.24%  "
$ws.Range("D38").Value = "'4.64"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.25%  "
$ws.Range("E39").Value = "  -0.34%  "
$ws.Range("E40").Value = "  +0.06%  "
$ws.Range("D41").Value = "'121.46"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.65%  "
$ws.Range("E42").Value = "  +0.53%  "
$ws.Range("D43").Value = "'21.60"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.68%  "
$ws.Range("D44").Value = "'0.0304"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.67%  "
$ws.Range("D45").Value = "'2.012.92"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.02%  "
$ws.Range("D46").Value = "'3.16"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.13%  "
$ws.Range("E47").Value = "  +3.99%  "
$ws.Range("E48").Value = "  -2.53%  "
$ws.Range("E49").Value = "  -1.10%  "
$ws.Range("E50").Value = "  -0.18%  "
$ws.Range("D51").Value = "'80.61"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.00%  "
